$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right
$ws.Range("A1").EntireColumn.Insert()

# Set the new column's header text
$ws.Range("A1").Value = "Owners"

# Match the column width used in the target layout (closest value the
# host's character-width -> pixel rounding can reach to 24.85546875)
$ws.Range("A1").EntireColumn.ColumnWidth = 24

# Update the active selection to match the target worksheet view
$ws.Range("C6").Select()
